$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 269
$wsExpo.Range("F4").Value = 2690
$wsExpo.Range("F6").Value = 574

# Sheet "全部类型" (all types) - same events duplicated, update column F accordingly
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 269
$wsAll.Range("F6").Value = 2690
$wsAll.Range("F8").Value = 574
